$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Resting Rate" values (columns D, G, J) for rows 4-6
$ws.Range("D4").Value = 6.7
$ws.Range("G4").Value = 8.8
$ws.Range("J4").Value = 10.5

$ws.Range("D5").Value = 1.9
$ws.Range("G5").Value = 1.8
$ws.Range("J5").Value = 3.1

$ws.Range("D6").Value = 8
$ws.Range("G6").Value = 5.2
$ws.Range("J6").Value = 6.4

# Row 7 only has the G "Resting Rate" column populated
$ws.Range("G7").Value = 0.4

# Update the selected cell/view state to match the new active cell
$ws.Range("G8").Select()
